# Update gh-pages output (data refresh) for 江西-漫展信息.xlsx
# Applies the same "want-to-go count" / price / venue / cover refresh to
# both the "展览" sheet and the combined "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 3 - 景德镇·江报国风动漫展
    $ws.Range("F3").Value = 1049

    # Row 4 - 南昌·晨啼漫拥二次元随机舞蹈派对·玺悦城场(免费活动)
    $ws.Range("F4").Value = 42

    # Row 6 - 江西·ShiningStaR动漫游戏文化节5th (new venue + updated cover)
    $ws.Range("D6").Value = "江西科技学院内 江西科技学院体育馆"
    $ws.Range("G6").Value = 52.1
    $ws.Range("I6").Value = "//i0.hdslb.com/bfs/openplatform/202403/p3TpZeAQ1709544877660.jpeg"

    # Row 10 - 南昌·AP动漫游戏 嘉年华内场票-小N&子音
    $ws.Range("F10").Value = 100

    # Row 11 - 南昌·CM01动漫游戏博览会
    $ws.Range("F11").Value = 771

    # Row 13 - 新余·文旅国漫嘉年华暨BM次元盛典 (updated cover)
    $ws.Range("F13").Value = 27
    $ws.Range("I13").Value = "//i2.hdslb.com/bfs/openplatform/202403/aXc6vPDP1709547191851.jpeg"

    # Row 14 - 赣州·第三届半夏动漫展
    $ws.Range("F14").Value = 206

    # Row 16 - 南昌·原X穹X崩only
    $ws.Range("F16").Value = 91

    # Row 17 - 南昌·第二届漫拥动漫嘉年华mini
    $ws.Range("F17").Value = 30
}

# Sheet-specific "want to go" counts that diverge slightly between the two sheets
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F6").Value = 2934
$wsExpo.Range("F8").Value = 1937

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 2936
$wsAll.Range("F8").Value = 1938
